$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.888.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.19%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.841.63"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.56%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.57%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.53"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.15%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4707"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.43%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3651"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.52%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07149"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.36%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9185"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.79%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.51"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.82%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07630"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.11%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.846.14"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.96%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.279"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.41%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.400"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.66%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.00"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.56%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.56%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008632"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.89%  "

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.45%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.922.51"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.24%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.013"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.06%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.61%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.926"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.22%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.75"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.02%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.19"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.008"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.21"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.62%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.855"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.40%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08824"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.16%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.224"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.43%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.171"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.41%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7434"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.21%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.474"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.04%  "

# Row 35
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.747"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.42%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.60%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01942"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.20%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05234"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.07%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.969"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.73%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5182"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.63%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.956"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.54%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1511"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.04%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.152"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.54%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.48"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.36%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4699"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.14%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.51%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.00"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.81%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.594"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.65%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.64"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.93%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06034"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8850"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.56%  "
